$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clients")

# Copy the formatting of the last data row down to the rows we are about to add,
# so new cells pick up the existing text-format + border style instead of a
# brand-new style entry.
$ws.Range("A4:B4").Copy()
$ws.Range("A5:B8").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Full refreshed client list (column A), "Yes" execution status (column B)
$clients = @("i2o-dev-su", "i2o-dev-cedar", "i2o-dev-dell", "i2o-dev-victrola", "i2o-dev-ecovacs", "i2o-dev-ausgold", "i2o-dev-jvc")

for ($i = 0; $i -lt $clients.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $clients[$i]
    $ws.Cells.Item($row, 2).Value = "Yes"
}

# Drop the now-unused "Last Run Date And Time" column entirely
$ws.Columns.Item(3).Delete()

# Switch the print setup to portrait orientation
$ws.PageSetup.Orientation = 1

$ws.Range("B2:B8").Select()
